$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# The log entry's Subject was re-selected from "Internal Medicine 1" to
# "Immuno&hema" for both checklist rows.
$ws.Range("B2").Value = "Immuno&hema"
$ws.Range("B3").Value = "Immuno&hema"

# Excel highlighted the edited Subject cells with a light-gray fill and
# centered the text (both horizontally and vertically) at a slightly
# smaller font size.
$rng = $ws.Range("B2:B3")
$rng.Interior.Color = 15790320
$rng.Interior.PatternColor = 15790320
$rng.Font.Size = 11
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$ws.Range("A1").Select()
